$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: was currency_conversion_rate / mdex:double / 6 / Currency Conversion Rate
# now po_currency_code / mdex:string / 1 / PO Currency
$ws.Range("C2").Value = "mdex:string"
$ws.Range("E2").Value = "PO Currency"

# Add new row 3 with the values that used to live on row 2
$ws.Range("A3").Value = 204
$ws.Range("E3").Value = "Func PO Unit Price"
$ws.Range("C3").Value = "mdex:double"
$ws.Range("D3").Value = 6

$ws.Range("B2").Value = "po_currency_code"
$ws.Range("D2").Value = 1

$ws.Range("B3").Value = "func_item_cost"

$ws.Range("B3").Select()
